# Update countries & provincias Spain
# - Update "Datos actualizados..." timestamp
# - Insert fresh "Banglades" row (pushing Eslovenia..Hong Kong down one row)
# - Insert fresh "Timor Oriental" row (pushing Santa Sede..Sahara Occidental down one row)
# - Update Austria (row 20) and Sri Lanka (row 117) figures with new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 10:52"

# --- Austria (row 20) ---
$ws.Range("D20").Value = 8098
$ws.Range("E20").Value = 5774
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 393

# --- New "Banglades" row inserted before Eslovenia; rows 70-76 shift to 71-77 ---
$ws.Range("A70").Value = "Banglades"
$ws.Range("B70").Value = 1231
$ws.Range("C70").Value = 219
$ws.Range("D70").Value = 49
$ws.Range("E70").Value = 1132
$ws.Range("F70").Value = 1
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 50

$ws.Range("A71").Value = "Eslovenia"
$ws.Range("B71").Value = 1220
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 152
$ws.Range("E71").Value = 1012
$ws.Range("F71").Value = 35
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 56

$ws.Range("A72").Value = "Uzbekistan"
$ws.Range("B72").Value = 1214
$ws.Range("C72").Value = 49
$ws.Range("D72").Value = 99
$ws.Range("E72").Value = 1111
$ws.Range("F72").Value = 8
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 4

$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 1197
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 351
$ws.Range("E73").Value = 833
$ws.Range("F73").Value = 25
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 13

$ws.Range("A74").Value = "Armenia"
$ws.Range("B74").Value = 1111
$ws.Range("C74").Value = 44
$ws.Range("D74").Value = 297
$ws.Range("E74").Value = 797
$ws.Range("F74").Value = 30
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 17

$ws.Range("A75").Value = "Lituania"
$ws.Range("B75").Value = 1091
$ws.Range("C75").Value = 21
$ws.Range("D75").Value = 138
$ws.Range("E75").Value = 924
$ws.Range("F75").Value = 14
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 29

$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 1083
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 236
$ws.Range("E76").Value = 807
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 40

$ws.Range("A77").Value = "Hong Kong"
$ws.Range("B77").Value = 1017
$ws.Range("C77").Value = 4
$ws.Range("D77").Value = 459
$ws.Range("E77").Value = 554
$ws.Range("F77").Value = 10
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 4

# --- Sri Lanka (row 117) ---
$ws.Range("D117").Value = 63
$ws.Range("E117").Value = 163

# --- New "Timor Oriental" row inserted before Santa Sede; rows 202-204 shift to 203-205 ---
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("B202").Value = 8
$ws.Range("C202").Value = 2
$ws.Range("D202").Value = 1
$ws.Range("E202").Value = 7
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Santa Sede"
$ws.Range("B203").Value = 8
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 6
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

$ws.Range("A204").Value = "Mauritania"
$ws.Range("B204").Value = 7
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 2
$ws.Range("E204").Value = 4
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 1

$ws.Range("A205").Value = "Sahara Occidental"
$ws.Range("B205").Value = 6
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 6
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0
